# Bold the "${company.denomination}" placeholder in the "Nous soussigné(s)
# ${company.denomination} , paierai contre le présent billet" paragraph.
# The trailing space right after the closing "}" must stay non-bold, so we
# search for the exact "${company.denomination}" substring (no trailing
# space) and only bold that found range.

$d = $word.ActiveDocument

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Nous soussign*") {
        $target = $para
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $found = $r.Find.Execute('${company.denomination}', $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.Font.Bold = 1
    }
}
